$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (IP) onto the new header cells, then set their text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..36: I = 1 (constant), J = same value as H (IP)
for ($r = 2; $r -le 36; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
